$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.717.26'
$ws.Cells.Item(2, 5).Value = '  -0.74%  '

$ws.Cells.Item(3, 4).Value = '1.848.26'
$ws.Cells.Item(3, 5).Value = '  -0.88%  '

$ws.Cells.Item(4, 5).Value = '  +0.13%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '315.04'
$ws.Cells.Item(5, 5).Value = '  -0.80%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.003'
$ws.Cells.Item(6, 5).Value = '  +0.16%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4266'
$ws.Cells.Item(7, 5).Value = '  -2.09%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3663'
$ws.Cells.Item(8, 5).Value = '  -1.69%  '

$ws.Cells.Item(9, 5).Value = '  +0.60%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.07319'
$ws.Cells.Item(10, 5).Value = '  -1.99%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.8921'
$ws.Cells.Item(11, 5).Value = '  -4.55%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '20.85'
$ws.Cells.Item(12, 5).Value = '  -2.16%  '

$ws.Cells.Item(13, 4).Value = '1.922.26'
$ws.Cells.Item(13, 5).Value = '  +1.22%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.579'
$ws.Cells.Item(14, 5).Value = '  -2.30%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '5.350'
$ws.Cells.Item(15, 5).Value = '  -1.72%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.06920'
$ws.Cells.Item(16, 5).Value = '  +0.79%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '1.005'
$ws.Cells.Item(17, 5).Value = '  +0.21%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '79.19'
$ws.Cells.Item(18, 5).Value = '  -2.94%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.000008901'
$ws.Cells.Item(19, 5).Value = '  -1.78%  '

$ws.Cells.Item(20, 5).Value = '  +0.19%  '

$ws.Cells.Item(21, 5).Value = '  -2.46%  '

$ws.Cells.Item(22, 4).Value = '27.729.11'
$ws.Cells.Item(22, 5).Value = '  -0.68%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '4.991'
$ws.Cells.Item(23, 5).Value = '  -2.52%  '

$ws.Cells.Item(24, 5).Value = '  -3.70%  '

$ws.Cells.Item(25, 4).Value = '2.048.64'
$ws.Cells.Item(25, 5).Value = '  -3.36%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.965'
$ws.Cells.Item(26, 5).Value = '  -1.96%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '153.82'
$ws.Cells.Item(27, 5).Value = '  -0.44%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '18.93'
$ws.Cells.Item(28, 5).Value = '  +2.77%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '120.91'
$ws.Cells.Item(29, 5).Value = '  +6.67%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '5.240'
$ws.Cells.Item(30, 5).Value = '  -4.11%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.902'
$ws.Cells.Item(31, 5).Value = '  +10.96%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.08942'
$ws.Cells.Item(32, 5).Value = '  -0.84%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.7705'
$ws.Cells.Item(33, 5).Value = '  -6.14%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.585'
$ws.Cells.Item(34, 5).Value = '  -4.82%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.968'
$ws.Cells.Item(35, 5).Value = '  -0.14%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.105'
$ws.Cells.Item(36, 5).Value = '  -6.07%  '

$ws.Cells.Item(37, 5).Value = '  -0.03%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.05397'
$ws.Cells.Item(38, 5).Value = '  -1.79%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.096'
$ws.Cells.Item(39, 5).Value = '  -2.24%  '

$ws.Cells.Item(40, 5).Value = '  -1.02%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.813'
$ws.Cells.Item(41, 5).Value = '  -5.38%  '

$ws.Cells.Item(42, 2).Value = 'FraxShare'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '6.926'
$ws.Cells.Item(42, 5).Value = '  -1.62%  '

$ws.Cells.Item(43, 2).Value = 'TheSandbox'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.5121'
$ws.Cells.Item(43, 5).Value = '  -2.71%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.1663'
$ws.Cells.Item(44, 5).Value = '  -2.41%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '8.283'
$ws.Cells.Item(45, 5).Value = '  -5.54%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.06599'
$ws.Cells.Item(46, 5).Value = '  -2.32%  '

$ws.Cells.Item(47, 2).Value = 'EnergySwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '10.46'
$ws.Cells.Item(47, 5).Value = '  -1.51%  '

$ws.Cells.Item(48, 2).Value = 'Decentraland'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.4761'
$ws.Cells.Item(48, 5).Value = '  -2.85%  '

$ws.Cells.Item(49, 5).Value = '  -2.28%  '

$ws.Cells.Item(50, 5).Value = '  +0.18%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.636'
$ws.Cells.Item(51, 5).Value = '  -2.46%  '
